$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing headers
$ws.Range("C1").Value = "Qurator"
$ws.Range("D1").Value = "Novischförman"

# Add new header column E, using the same style as the data cells (e.g. D2)
$ws.Range("D2").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E1").Value = "PQE"

$ws.Range("E1").Select()
